$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "64.325.79"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.484.30"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "586.04"
$ws.Range("D6").Value = "134.02"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D7").Value = "3.483.32"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "7.18"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "4.074.75"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "3.485.55"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "64.345.82"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "24.98"
$ws.Range("E18").Value = "  -10.19%  "
$ws.Range("D19").Value = "9.96"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "13.72"
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").Value = "385.08"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "0.565"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "3.620.59"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "5.72"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000113"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "1.55"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "7.41"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "8.18"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.500.58"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.146"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "23.35"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "5.27"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.82"
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "1.53"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "162.87"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0777"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.802"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "25.49"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "41.76"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.66"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.467.31"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "6.74"
$ws.Range("E51").Value = "  -2.06%  "
